# Log today's walk: bump the running February distance total (column G)
# on the data sheet, then flip back to the chart tab to see the new point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # "Sheet1" is the active tab in this workbook
$ws.Range("G2").Value = 64.5   # add today's distance to the Feb running total
$ws.Range("G3").Select()       # leave the selection where entry would continue

$wb.Worksheets.Item("Chart1").Activate()
